# Update guide with new Message system
#
# The "Send Messages" leg of the diagram (an elbow connector plus its two
# labels) is removed. One of the two labels ("Accept input") stays on the
# slide; the connector that routed into it ("Connector: Elbow 110") and the
# second label ("TextBox 112" / "Send Messages") are deleted together with
# it, per the commit.
#
# Also refresh any cached "datetimeFigureOut" date fields (present on
# master/layout date placeholders instantiated on slides) from 26/5/2020 to
# 26/10/2020, as seen throughout the canonical diff.

$p = $ppt.ActivePresentation

$namesToDelete = @("Connector: Elbow 110", "TextBox 112")

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)

    # Remove the obsolete connector + "Send Messages" label.
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $shp = $s.Shapes.Item($i)
        if ($namesToDelete -contains $shp.Name) {
            $shp.Delete()
        }
    }

    # Refresh any cached date fields left on the slide.
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "26/5/2020") {
                $tr.Text = "26/10/2020"
            }
        }
    }
}
